$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '66.548.64'
$ws.Cells.Item(2, 5).Value = '  +1.14%  '

$ws.Cells.Item(3, 4).Value = '3.307.78'
$ws.Cells.Item(3, 5).Value = '  +0.43%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).Value = '591.09'
$ws.Cells.Item(5, 5).Value = '  +3.17%  '

$ws.Cells.Item(6, 4).Value = '180.99'
$ws.Cells.Item(6, 5).Value = '  +1.55%  '

$ws.Cells.Item(7, 4).Value = '0.637'
$ws.Cells.Item(7, 5).Value = '  +1.13%  '

$ws.Cells.Item(8, 5).Value = '  -0.02%  '

$ws.Cells.Item(9, 4).Value = '3.304.33'
$ws.Cells.Item(9, 5).Value = '  +0.40%  '

$ws.Cells.Item(10, 5).Value = '  +0.33%  '

$ws.Cells.Item(11, 4).Value = '6.87'
$ws.Cells.Item(11, 5).Value = '  +3.04%  '

$ws.Cells.Item(12, 5).Value = '  +0.47%  '

$ws.Cells.Item(13, 4).Value = '3.881.27'
$ws.Cells.Item(13, 5).Value = '  +0.37%  '

$ws.Cells.Item(14, 5).Value = '  -2.13%  '

$ws.Cells.Item(15, 4).Value = '66.562.66'
$ws.Cells.Item(15, 5).Value = '  +0.90%  '

$ws.Cells.Item(16, 4).Value = '26.72'
$ws.Cells.Item(16, 5).Value = '  +0.43%  '

$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '3.320.86'
$ws.Cells.Item(17, 5).Value = '  +0.12%  '

$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).Value = '0.0000163'
$ws.Cells.Item(18, 5).Value = '  -0.01%  '

$ws.Cells.Item(19, 4).Value = '429.76'
$ws.Cells.Item(19, 5).Value = '  -1.41%  '

$ws.Cells.Item(20, 4).Value = '5.48'
$ws.Cells.Item(20, 5).Value = '  -2.04%  '

$ws.Cells.Item(21, 4).Value = '13.05'
$ws.Cells.Item(21, 5).Value = '  -1.98%  '

$ws.Cells.Item(22, 4).Value = '7.30'
$ws.Cells.Item(22, 5).Value = '  -1.75%  '

$ws.Cells.Item(23, 5).Value = '  +0.15%  '

$ws.Cells.Item(24, 4).Value = '71.46'
$ws.Cells.Item(24, 5).Value = '  -1.31%  '

$ws.Cells.Item(25, 5).Value = '  +1.00%  '

$ws.Cells.Item(26, 4).Value = '3.448.75'

$ws.Cells.Item(27, 5).Value = '  +0.13%  '

$ws.Cells.Item(28, 5).Value = '  +5.84%  '

$ws.Cells.Item(29, 5).Value = '  +0.65%  '

$ws.Cells.Item(30, 4).Value = '9.21'
$ws.Cells.Item(30, 5).Value = '  +2.79%  '

$ws.Cells.Item(31, 5).Value = '  -0.06%  '

$ws.Cells.Item(32, 5).Value = '  -0.85%  '

$ws.Cells.Item(33, 4).Value = '22.35'
$ws.Cells.Item(33, 5).Value = '  -0.24%  '

$ws.Cells.Item(34, 5).Value = '  +0.09%  '

$ws.Cells.Item(35, 4).Value = '5.18'
$ws.Cells.Item(35, 5).Value = '  +0.63%  '

$ws.Cells.Item(36, 4).Value = '6.58'
$ws.Cells.Item(36, 5).Value = '  -0.91%  '

$ws.Cells.Item(37, 5).Value = '  -0.82%  '

$ws.Cells.Item(38, 4).Value = '159.13'
$ws.Cells.Item(38, 5).Value = '  +1.30%  '

$ws.Cells.Item(39, 5).Value = '  -1.59%  '

$ws.Cells.Item(40, 4).Value = '2.869.73'
$ws.Cells.Item(40, 5).Value = '  +3.11%  '

$ws.Cells.Item(41, 5).Value = '  +0.23%  '

$ws.Cells.Item(42, 4).Value = '26.36'
$ws.Cells.Item(42, 5).Value = '  -2.64%  '

$ws.Cells.Item(43, 4).Value = '4.33'
$ws.Cells.Item(43, 5).Value = '  -0.57%  '

$ws.Cells.Item(44, 4).Value = '0.751'
$ws.Cells.Item(44, 5).Value = '  -4.15%  '

$ws.Cells.Item(45, 4).Value = '39.78'
$ws.Cells.Item(45, 5).Value = '  -1.40%  '

$ws.Cells.Item(46, 2).Value = 'RenderToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(46, 4).Value = '5.94'
$ws.Cells.Item(46, 5).Value = '  -2.59%  '

$ws.Cells.Item(47, 2).Value = 'Hedera'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(47, 4).Value = '0.0653'
$ws.Cells.Item(47, 5).Value = '  -0.86%  '

$ws.Cells.Item(48, 4).Value = '2.32'
$ws.Cells.Item(48, 5).Value = '  +0.90%  '

$ws.Cells.Item(49, 4).Value = '313.40'
$ws.Cells.Item(49, 5).Value = '  -2.88%  '

$ws.Cells.Item(50, 4).Value = '22.90'
$ws.Cells.Item(50, 5).Value = '  -2.74%  '

$ws.Cells.Item(51, 4).Value = '0.0271'
$ws.Cells.Item(51, 5).Value = '  +0.02%  '
